# Daily attendance processing - 2025-10-15 08:28:54
# Updates the CNS session-analysis sheet with refreshed attendance data:
# new session date/recorder/headcount rows, re-ordered recorder-email
# lists, and recalculated attendance percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ANATOMY, session 1) -------------------------------------------
# Session date moved out a couple of weeks; recorders + headcount updated.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "08/10/2025"
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H2").Value = "157/221"

# --- Row 3 (ANATOMY, session 2) -------------------------------------------
$ws.Range("G3").Value = "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# --- Row 4 (ANATOMY, session 3) -------------------------------------------
$ws.Range("G4").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# --- Row 5 (ANATOMY, session 4) -------------------------------------------
$ws.Range("G5").Value = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, hananragab@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("H5").Value = "78/221"

# --- Class statistics ------------------------------------------------------
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "47.0%"

# --- Row 12 ------------------------------------------------------------
$ws.Range("G12").Value = "salma.elgendy.std@med.asu.edu.eg, System"

# --- Row 13 ------------------------------------------------------------
$ws.Range("G13").Value = "mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# --- Attendance percentages --------------------------------------------
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "54.4%"
$ws.Range("S16").NumberFormat = "@"
$ws.Range("S16").Value = "38.6%"

# --- Row 24 ------------------------------------------------------------
$ws.Range("G24").Value = "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# --- Row 25 ------------------------------------------------------------
$ws.Range("G25").Value = "marina_atef@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# --- Row 31 ------------------------------------------------------------
$ws.Range("G31").Value = "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# --- Row 32 ------------------------------------------------------------
$ws.Range("G32").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# --- Row 33 ------------------------------------------------------------
$ws.Range("G33").Value = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, hananragab@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("H33").Value = "91/246"

# --- Row 40 ------------------------------------------------------------
$ws.Range("G40").Value = "salma.elgendy.std@med.asu.edu.eg, System"

# --- Row 41 ------------------------------------------------------------
$ws.Range("G41").Value = "mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# --- Row 52 ------------------------------------------------------------
$ws.Range("G52").Value = "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# --- Row 53 ------------------------------------------------------------
$ws.Range("G53").Value = "marina_atef@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
